$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4 for "Brighton Beach" (pushes the existing
# Brighton restaurant row, and everything below, down by one).
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "Brighton"
$ws.Range("B4").Value = "Brighton Beach"
$ws.Range("C4").Value = "26/12/20 12pm - 1pm"
$ws.Range("D4").Value = "Case attended beach"

# Insert a new row at row 9 for "Cheltenham" (after the shift above,
# row 9 is still the old "Hampton" row; push it and everything below
# down by one more).
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "Cheltenham"
$ws.Range("B9").Value = "Two Bob Snob, 256 Charman Road"
$ws.Range("C9").Value = "22/12/2020 1pm - 2pm"
$ws.Range("D9").Value = "Case attended Venue"

# Update the exposure period for the Fonda Mexican Flinders Lane row
# (originally row 22, now row 24 after the two insertions above).
$ws.Range("C24").Value = "29/12/20 6:00pm-7:30pm"
